$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.46628438110641
$ws.Range("D2").Value = 9.071548159512837
$ws.Range("E2").Value = 14.18880604075213
$ws.Range("F2").Value = 34.99454784501328
$ws.Range("G2").Value = 37.62470869637943
$ws.Range("H2").Value = 16.67422024450506
$ws.Range("J2").Value = 10.63970615628701
$ws.Range("K2").Value = 15.76039444621195
$ws.Range("L2").Value = 9.756038640054662
$ws.Range("M2").Value = 18.43099064586145
$ws.Range("O2").Value = 26.52540931578303
$ws.Range("C3").Value = 13.45300446102775
$ws.Range("D3").Value = 9.049227845633121
$ws.Range("E3").Value = 14.21733002514112
$ws.Range("F3").Value = 35.15162831205756
$ws.Range("G3").Value = 37.85303048725027
$ws.Range("H3").Value = 16.74862410732115
$ws.Range("J3").Value = 10.66867351873471
$ws.Range("K3").Value = 15.23686830112957
$ws.Range("L3").Value = 9.773625342960077
$ws.Range("M3").Value = 18.22158603789687
$ws.Range("O3").Value = 26.66434367173343
$ws.Range("C4").Value = 13.4475303277871
$ws.Range("D4").Value = 9.036519112148378
$ws.Range("E4").Value = 14.23672885954186
$ws.Range("F4").Value = 35.25653683882179
$ws.Range("G4").Value = 38.00509060651048
$ws.Range("H4").Value = 16.79721048790941
$ws.Range("J4").Value = 10.68746015419034
$ws.Range("K4").Value = 14.90603542088003
$ws.Range("L4").Value = 9.785049994002797
$ws.Range("M4").Value = 18.09302545641983
$ws.Range("O4").Value = 26.75562000332408
$ws.Range("C5").Value = 13.44597604458243
$ws.Range("D5").Value = 9.031593890861425
$ws.Range("E5").Value = 14.24510835416974
$ws.Range("F5").Value = 35.30141068026172
$ws.Range("G5").Value = 38.07002928406212
$ws.Range("H5").Value = 16.81774000851946
$ws.Range("J5").Value = 10.69536807009267
$ws.Range("K5").Value = 14.76903356329109
$ws.Range("L5").Value = 9.789863556348505
$ws.Range("M5").Value = 18.04068524552943
$ws.Range("O5").Value = 26.79431558957212
$ws.Range("C6").Value = 13.44575888200023
$ws.Range("D6").Value = 9.030791474926106
$ws.Range("E6").Value = 14.24652841893045
$ws.Range("F6").Value = 35.30899001937907
$ws.Range("G6").Value = 38.0809915175923
$ws.Range("H6").Value = 16.82119303335986
$ws.Range("J6").Value = 10.69669642434114
$ws.Range("K6").Value = 14.74615800060036
$ws.Range("L6").Value = 9.790672395457905
$ws.Range("M6").Value = 18.03199850001038
$ws.Range("O6").Value = 26.80083150062696
$ws.Range("C7").Value = 13.44750662404189
$ws.Range("D7").Value = 9.036451657561345
$ws.Range("E7").Value = 14.23683994746175
$ws.Range("F7").Value = 35.25713343439151
$ws.Range("G7").Value = 38.00595437042987
$ws.Range("H7").Value = 16.79748439918662
$ws.Range("J7").Value = 10.6875657811113
$ws.Range("K7").Value = 14.90419636825009
$ws.Range("L7").Value = 9.785114271393033
$ws.Range("M7").Value = 18.09231931837759
$ws.Range("O7").Value = 26.75613579507586
$ws.Range("C8").Value = 13.46115097087259
$ws.Range("D8").Value = 9.063647376504173
$ws.Range("E8").Value = 14.19825007653617
$ws.Range("F8").Value = 35.04695041109242
$ws.Range("G8").Value = 37.70096342839715
$ws.Range("H8").Value = 16.69927261290589
$ws.Range("J8").Value = 10.64948680542218
$ws.Range("K8").Value = 15.58193077206089
$ws.Range("L8").Value = 9.761972793428757
$ws.Range("M8").Value = 18.35881221916103
$ws.Range("O8").Value = 26.57207389121162
$ws.Range("C9").Value = 13.50903766250231
$ws.Range("D9").Value = 9.124732994581537
$ws.Range("E9").Value = 14.13751770170549
$ws.Range("F9").Value = 34.70212007193291
$ws.Range("G9").Value = 37.19758412673317
$ws.Range("H9").Value = 16.52968744353594
$ws.Range("J9").Value = 10.58272512207215
$ws.Range("K9").Value = 16.82961510129329
$ws.Range("L9").Value = 9.721542884690836
$ws.Range("M9").Value = 18.87952650833104
$ws.Range("O9").Value = 26.25856589438926
$ws.Range("C10").Value = 13.55689749790904
$ws.Range("D10").Value = 9.174132839401754
$ws.Range("E10").Value = 14.10198680205436
$ws.Range("F10").Value = 34.49010538811007
$ws.Range("G10").Value = 36.88620941375243
$ws.Range("H10").Value = 16.41909274417965
$ws.Range("J10").Value = 10.53845902923895
$ws.Range("K10").Value = 17.68858955625462
$ws.Range("L10").Value = 9.69483007005649
$ws.Range("M10").Value = 19.25830661984395
$ws.Range("O10").Value = 26.05724315261213
$ws.Range("C11").Value = 13.58137033788832
$ws.Range("D11").Value = 9.197540201552734
$ws.Range("E11").Value = 14.0877917534724
$ws.Range("F11").Value = 34.40268871821092
$ws.Range("G11").Value = 36.75740641697566
$ws.Range("H11").Value = 16.37181511458727
$ws.Range("J11").Value = 10.51935170247666
$ws.Range("K11").Value = 18.06541330235375
$ws.Range("L11").Value = 9.683321492629698
$ws.Range("M11").Value = 19.42925167729902
$ws.Range("O11").Value = 25.97197724282431
$ws.Range("C12").Value = 13.59102079720131
$ws.Range("D12").Value = 9.206534143731531
$ws.Range("E12").Value = 14.0826990338667
$ws.Range("F12").Value = 34.37088969940159
$ws.Range("G12").Value = 36.71049192845204
$ws.Range("H12").Value = 16.35434811825755
$ws.Range("J12").Value = 10.51226368569474
$ws.Range("K12").Value = 18.20600367164208
$ws.Range("L12").Value = 9.679055555597742
$ws.Range("M12").Value = 19.49374312510381
$ws.Range("O12").Value = 25.94059967869844
$ws.Range("C13").Value = 13.58892544743433
$ws.Range("D13").Value = 9.204591423292662
$ws.Range("E13").Value = 14.08378327816587
$ws.Range("F13").Value = 34.37768011779377
$ws.Range("G13").Value = 36.72051284394611
$ws.Range("H13").Value = 16.3580905569139
$ws.Range("J13").Value = 10.51378366466808
$ws.Range("K13").Value = 18.17582011304624
$ws.Range("L13").Value = 9.679970210831842
$ws.Range("M13").Value = 19.47986520284342
$ws.Range("O13").Value = 25.94731685153782
$ws.Range("C14").Value = 13.5821566404477
$ws.Range("D14").Value = 9.198277557183909
$ws.Range("E14").Value = 14.08736711017883
$ws.Range("F14").Value = 34.40004643412847
$ws.Range("G14").Value = 36.7535093703952
$ws.Range("H14").Value = 16.37036935603133
$ws.Range("J14").Value = 10.51876561371944
$ws.Range("K14").Value = 18.07702247528467
$ws.Range("L14").Value = 9.682968687548565
$ws.Range("M14").Value = 19.4345625293336
$ws.Range("O14").Value = 25.9693775240257
$ws.Range("C15").Value = 13.57806028042514
$ws.Range("D15").Value = 9.194426939771985
$ws.Range("E15").Value = 14.08959910668612
$ws.Range("F15").Value = 34.4139163946774
$ws.Range("G15").Value = 36.77396338003156
$ws.Range("H15").Value = 16.37794725634128
$ws.Range("J15").Value = 10.52183639621118
$ws.Range("K15").Value = 18.01622917791521
$ws.Range("L15").Value = 9.684817325465929
$ws.Range("M15").Value = 19.40678051893936
$ws.Range("O15").Value = 25.98300900762385
$ws.Range("C16").Value = 13.55535205157967
$ws.Range("D16").Value = 9.172621558551059
$ws.Range("E16").Value = 14.10295405419919
$ws.Range("F16").Value = 34.49600043132544
$ws.Range("G16").Value = 36.89488658873182
$ws.Range("H16").Value = 16.42224346439231
$ws.Range("J16").Value = 10.53972840921552
$ws.Range("K16").Value = 17.66367398628065
$ws.Range("L16").Value = 9.695595093054875
$ws.Range("M16").Value = 19.24710358829509
$ws.Range("O16").Value = 26.06294279723372
$ws.Range("C17").Value = 13.54210948984847
$ws.Range("D17").Value = 9.159481095535957
$ws.Range("E17").Value = 14.11165070996157
$ws.Range("F17").Value = 34.54867288542741
$ws.Range("G17").Value = 36.97236812711339
$ws.Range("H17").Value = 16.4501943934075
$ws.Range("J17").Value = 10.55096788809237
$ws.Range("K17").Value = 17.44374908953854
$ws.Range("L17").Value = 9.702371374358226
$ws.Range("M17").Value = 19.14876482734026
$ws.Range("O17").Value = 26.11359918449931
$ws.Range("C18").Value = 13.53474730962351
$ws.Range("D18").Value = 9.152011339117552
$ws.Range("E18").Value = 14.11683807128764
$ws.Range("F18").Value = 34.57981831026071
$ws.Range("G18").Value = 37.01814148038617
$ws.Range("H18").Value = 16.4665564581242
$ws.Range("J18").Value = 10.55752947580895
$ws.Range("K18").Value = 17.31594838756062
$ws.Range("L18").Value = 9.706329477293664
$ws.Range("M18").Value = 19.09207725616294
$ws.Range("O18").Value = 26.14332981671899
$ws.Range("C19").Value = 13.53229848243877
$ws.Range("D19").Value = 9.149497503215846
$ws.Range("E19").Value = 14.11862625484259
$ws.Range("F19").Value = 34.59050937696967
$ws.Range("G19").Value = 37.03384661544752
$ws.Range("H19").Value = 16.47214539690711
$ws.Range("J19").Value = 10.5597677820126
$ws.Range("K19").Value = 17.27245631346591
$ws.Range("L19").Value = 9.707680036765575
$ws.Range("M19").Value = 19.07286363472604
$ws.Range("O19").Value = 26.1534981169041
$ws.Range("C20").Value = 13.54349287094023
$ws.Range("D20").Value = 9.160870813870858
$ws.Range("E20").Value = 14.11070576348585
$ws.Range("F20").Value = 34.54297784144529
$ws.Range("G20").Value = 36.96399496071494
$ws.Range("H20").Value = 16.447189428838
$ws.Range("J20").Value = 10.54976139767337
$ws.Range("K20").Value = 17.46729638371764
$ws.Range("L20").Value = 9.701643762158572
$ws.Range("M20").Value = 19.15924650152863
$ws.Range("O20").Value = 26.10814518813892
$ws.Range("C21").Value = 13.58413444957027
$ws.Range("D21").Value = 9.200128597989107
$ws.Range("E21").Value = 14.0863067843295
$ws.Range("F21").Value = 34.39344147893259
$ws.Range("G21").Value = 36.74376689122776
$ws.Range("H21").Value = 16.36675094195279
$ws.Range("J21").Value = 10.51729829527318
$ws.Range("K21").Value = 18.10609960484394
$ws.Range("L21").Value = 9.682085465402682
$ws.Range("M21").Value = 19.4478759409878
$ws.Range("O21").Value = 25.96287303151923
$ws.Range("C22").Value = 13.61292646249826
$ws.Range("D22").Value = 9.226541914342127
$ws.Range("E22").Value = 14.07200786795617
$ws.Range("F22").Value = 34.30331262970263
$ws.Range("G22").Value = 36.61068532565086
$ws.Range("H22").Value = 16.31672117699317
$ws.Range("J22").Value = 10.49694135798031
$ws.Range("K22").Value = 18.51128844941769
$ws.Range("L22").Value = 9.669839709107555
$ws.Range("M22").Value = 19.63508264648548
$ws.Range("O22").Value = 25.87323921737686
$ws.Range("C23").Value = 13.59735740338392
$ws.Range("D23").Value = 9.212376934141526
$ws.Range("E23").Value = 14.07948888092068
$ws.Range("F23").Value = 34.35071879890797
$ws.Range("G23").Value = 36.68071617088053
$ws.Range("H23").Value = 16.34319046955143
$ws.Range("J23").Value = 10.50772776179291
$ws.Range("K23").Value = 18.29618802552094
$ws.Range("L23").Value = 9.676326515049588
$ws.Range("M23").Value = 19.53531244504606
$ws.Range("O23").Value = 25.92059178699168
$ws.Range("C24").Value = 13.54286666167358
$ws.Range("D24").Value = 9.160242257529294
$ws.Range("E24").Value = 14.1111323895995
$ws.Range("F24").Value = 34.54554988199591
$ws.Range("G24").Value = 36.96777664452385
$ws.Range("H24").Value = 16.44854706154641
$ws.Range("J24").Value = 10.55030654088676
$ws.Range("K24").Value = 17.45665488945918
$ws.Range("L24").Value = 9.701972521456788
$ws.Range("M24").Value = 19.15450820443165
$ws.Range("O24").Value = 26.11060904741246
$ws.Range("C25").Value = 13.49384018807124
$ws.Range("D25").Value = 9.107398242331136
$ws.Range("E25").Value = 14.15234970322206
$ws.Range("F25").Value = 34.78816760392402
$ws.Range("G25").Value = 37.32354814218614
$ws.Range("H25").Value = 16.5731044539899
$ws.Range("J25").Value = 10.59994305570512
$ws.Range("K25").Value = 16.50168830754424
$ws.Range("L25").Value = 9.731953058383189
$ws.Range("M25").Value = 18.73913552017863
$ws.Range("O25").Value = 26.33829037251878
